$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# All edited cells are forced to Text format ("@") before assignment so that
# numeric-looking price strings (e.g. "610.59", "32.00", "0.0691") are preserved
# exactly as text instead of being auto-converted to floating point numbers,
# matching the original workbook formatting (inline/shared text strings).

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "89.477.19"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.34%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.022.36"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -3.38%  "

# Row 4
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.10%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "210.19"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -1.97%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "610.59"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -3.94%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.362"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -7.80%  "

# Row 8
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +14.13%  "

# Row 9
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.09%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "3.019.14"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -3.45%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.665"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +19.43%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.188"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +5.24%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000239"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -4.52%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.31"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -0.39%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "89.224.03"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +0.25%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "32.00"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -0.49%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.586.42"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -3.22%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.039.30"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -3.22%  "

# Row 19
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -2.36%  "

# Row 20
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -2.42%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.33"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.99%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "422.12"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.72%  "

# Row 23
$ws.Range("B23").NumberFormat = "@"
$ws.Range("B23").Value = "Uniswap"
$ws.Range("C23").NumberFormat = "@"
$ws.Range("C23").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.22"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -1.41%  "

# Row 24
$ws.Range("B24").NumberFormat = "@"
$ws.Range("B24").Value = "Polkadot"
$ws.Range("C24").NumberFormat = "@"
$ws.Range("C24").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.01"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +2.32%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "5.32"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -2.56%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "82.17"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +0.37%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.50"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +0.00%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.192.71"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -2.89%  "

# Row 29
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.02%  "

# Row 30
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +3.65%  "

# Row 31
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +2.60%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.24"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +0.63%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.75"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -5.78%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "498.36"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -1.06%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.61"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -5.64%  "

# Row 36
$ws.Range("B36").NumberFormat = "@"
$ws.Range("B36").Value = "PancakeSwap"
$ws.Range("C36").NumberFormat = "@"
$ws.Range("C36").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.80"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -1.84%  "

# Row 37
$ws.Range("B37").NumberFormat = "@"
$ws.Range("B37").Value = "EthereumClassic"
$ws.Range("C37").NumberFormat = "@"
$ws.Range("C37").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "22.62"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +3.11%  "

# Row 38
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -3.19%  "

# Row 39
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -10.55%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "22.24"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +0.03%  "

# Row 41
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +0.04%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.357"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -1.52%  "

# Row 44
$ws.Range("B44").NumberFormat = "@"
$ws.Range("B44").Value = "Stellar"
$ws.Range("C44").NumberFormat = "@"
$ws.Range("C44").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.135"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +3.37%  "

# Row 45
$ws.Range("B45").NumberFormat = "@"
$ws.Range("B45").Value = "Stacks"
$ws.Range("C45").NumberFormat = "@"
$ws.Range("C45").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.81"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -2.75%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "144.55"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -0.22%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "43.42"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -0.74%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0691"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +6.59%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "4.15"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +5.72%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "160.45"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -1.73%  "

# Row 51
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +2.02%  "
